$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 ("title" slide) - contact-info subtitle box
#   * paragraph 1: "hwdong"                -> "Youtube: hwdong"
#   * paragraph 2: "B站 和微博： hw-dong"   -> "B站： hw-dong"
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$contactBox = $s1.Shapes.Item(3)
$tf1 = $contactBox.TextFrame
$tr1 = $tf1.TextRange

# Add the "Youtube: " lead-in before the existing "hwdong" run.
[void]$tr1.Paragraphs(1, 1).InsertBefore("Youtube: ")

# Trim "B站 和微博： hw-dong" down to "B站： hw-dong" (replace the
# "站 和微博： " run's text with "站： ").
$tr1b = $contactBox.TextFrame.TextRange
$tr1b.Paragraphs(2, 1).Characters(2, 7).Text = "站： "

# ---------------------------------------------------------------------------
# Slide 28 ("赋值运算符") - nudge the lower screenshot into place and give
# it a click-triggered "Appear" entrance animation.
# ---------------------------------------------------------------------------
$s28 = $p.Slides.Item(28)
$shot = $s28.Shapes.Item(3)
$shot.Left = 22.44700813293457
$shot.Top = 300.24420166015625

$mainSeq = $s28.TimeLine.MainSequence
$mainSeq.AddEffect($shot, 1, 0, 1) | Out-Null
